$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Insert a new column before column D, shifting D onward to the right
$ws.Columns("D").Insert()

# Set header for new column D
$ws.Cells.Item(1, 4).Value = "Technical Replicate"

# New column D should visually match its neighbor (old D, now E) width
$ws.Columns("D").ColumnWidth = 21.5
